$wb = $excel.ActiveWorkbook
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsAllTime = $wb.Worksheets.Item("All Time")

# Updated Suzie's Roth IRA dividend for October (row 12) on the Yearly sheet.
# Downstream formulas (G12 shared-sum, F15/G15 totals) recalculate automatically.
$wsYearly.Range("F12").Value = 13.08

# Link the 2016 row on the "All Time" sheet to the Yearly totals instead of
# static numbers, so it stays in sync with the updated dividend figures.
$wsAllTime.Range("F7").Formula = "=Yearly!D15"
$wsAllTime.Range("G7").Formula = "=Yearly!E15"
$wsAllTime.Range("H7").Formula = "=Yearly!F15"

# Restore the selections that were active on each sheet.
$wsYearly.Activate()
$wsYearly.Range("D15").Select()

$wsAllTime.Activate()
$wsAllTime.Range("L15").Select()
